# "Bugs varios OK (relacionados a las clases extended)"
#
# - Marks rows 14, 20 and 24 as done (red -> green "ESTADO" marker).
# - Fills in three new pending TODO rows (25-27) with text + marks
#   their ESTADO cell green as well.
# - Leaves the selection on A27 (matches the author's last-edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Green fill used by the "done" cells in column A (matches the existing
# xf that uses fillId=5 / rgb FF00FF00 already present in the workbook).
$doneColor = 65280  # RGB(0,255,0) as a BGR-packed OLE color

# Rows that flip from "pending" (red) to "done" (green).
$doneRows = 14, 20, 24, 25, 26, 27
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 1).Interior.Color = $doneColor
}

# New TODO items added under "IMPORTANTE" section.
$ws.Range("B25").Value = "Deshabilitar legajo en personaABM modificar"
$ws.Range("B26").Value = "UsuarioABM alta mejorar combobox"
$ws.Range("B27").Value = "Listar planes, materias, comisiones, cursos,  edit y elim crash"

# Leave the view scrolled down a bit further with A27 selected, as in
# the author's saved workbook state.
$ws.Range("A27").Select() | Out-Null
